# "cleaned defensive actions data"
# - Row 1 becomes the real (flattened/abbreviated) header row, filling in the
#   cells that used to be blank placeholders under merged header groups.
# - The old sub-header (row 2), the blank spacer (row 3), and the totals row
#   (row 20) are hidden rather than removed.
# - The header merges (H1:L1, M1:P1, Q1:S1) are no longer needed once every
#   cell in row 1 carries its own label, so they are unmerged.
# - A handful of Tkl% cells that were blank (0/0 => NaN) are cleaned to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the header cell groups now that row 1 gets individual labels.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Populate row 1 with the real column headers.
$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# Hide the now-redundant old sub-header row, the blank spacer row, and the
# aggregate totals row.
$ws.Rows(2).EntireRow.Hidden = $true
$ws.Rows(3).EntireRow.Hidden = $true
$ws.Rows(20).EntireRow.Hidden = $true

# Clean up blank Tkl% values (0 tackles / 0 tackles won => NaN) to 0.
$ws.Range("O4").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("O19").Value = 0
